$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric-looking text columns (Price / Hora) must stay as TEXT so that
# values such as "0.1410" or "0.001600" keep their exact original formatting
# instead of being reinterpreted as floating point numbers. We flip the cell
# format to Text first, write the value, then restore the default "Normal"
# style so no stray formatting is left behind on cells that originally had none. ---
$numCells = @('D2', 'G2', 'D3', 'G3', 'D4', 'G4', 'D5', 'G5', 'D6', 'G6', 'D7', 'G7', 'D8', 'G8', 'D9', 'G9', 'D10', 'G10', 'D11', 'G11', 'D12', 'G12', 'D13', 'G13', 'D14', 'G14', 'D15', 'G15', 'D16', 'G16', 'D17', 'G17', 'D18', 'G18', 'G19', 'D20', 'G20', 'D21', 'G21', 'D22', 'G22', 'G23', 'D24', 'G24', 'G25', 'G26', 'G27', 'G28', 'G29', 'G30', 'G31', 'G32', 'G33', 'G34', 'G35', 'G36', 'G37', 'G38', 'G39', 'D40', 'G40', 'D41', 'G41', 'D42', 'G42', 'D43', 'G43', 'D44', 'G44', 'D45', 'G45', 'G46', 'D47', 'G47', 'D48', 'G48', 'G49', 'G50', 'G51')
foreach ($c in $numCells) {
    $ws.Range($c).NumberFormat = "@"
}

# --- Write updated cell values ---
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('E10').Value = '9OneONE'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('E18').Value = '17CoinExTokenCET'

$ws.Range('D2').Value = '246.78'
$ws.Range('G2').Value = '23'
$ws.Range('D3').Value = '22.55'
$ws.Range('G3').Value = '23'
$ws.Range('D4').Value = '5.263'
$ws.Range('G4').Value = '23'
$ws.Range('D5').Value = '0.05687'
$ws.Range('G5').Value = '23'
$ws.Range('D6').Value = '3.418'
$ws.Range('G6').Value = '23'
$ws.Range('D7').Value = '6.285'
$ws.Range('G7').Value = '23'
$ws.Range('D8').Value = '0.8095'
$ws.Range('G8').Value = '23'
$ws.Range('D9').Value = '0.8731'
$ws.Range('G9').Value = '23'
$ws.Range('D10').Value = '0.01095'
$ws.Range('G10').Value = '23'
$ws.Range('D11').Value = '0.1410'
$ws.Range('G11').Value = '23'
$ws.Range('D12').Value = '0.07339'
$ws.Range('G12').Value = '23'
$ws.Range('D13').Value = '0.03049'
$ws.Range('G13').Value = '23'
$ws.Range('D14').Value = '0.03075'
$ws.Range('G14').Value = '23'
$ws.Range('D15').Value = '0.09386'
$ws.Range('G15').Value = '23'
$ws.Range('D16').Value = '3.869'
$ws.Range('G16').Value = '23'
$ws.Range('D17').Value = '0.001600'
$ws.Range('G17').Value = '23'
$ws.Range('D18').Value = '0.04790'
$ws.Range('G18').Value = '23'
$ws.Range('G19').Value = '23'
$ws.Range('D20').Value = '0.005015'
$ws.Range('G20').Value = '23'
$ws.Range('D21').Value = '0.0009964'
$ws.Range('G21').Value = '23'
$ws.Range('D22').Value = '0.0001501'
$ws.Range('G22').Value = '23'
$ws.Range('G23').Value = '23'
$ws.Range('D24').Value = '2.189'
$ws.Range('G24').Value = '23'
$ws.Range('G25').Value = '23'
$ws.Range('G26').Value = '23'
$ws.Range('G27').Value = '23'
$ws.Range('G28').Value = '23'
$ws.Range('G29').Value = '23'
$ws.Range('G30').Value = '23'
$ws.Range('G31').Value = '23'
$ws.Range('G32').Value = '23'
$ws.Range('G33').Value = '23'
$ws.Range('G34').Value = '23'
$ws.Range('G35').Value = '23'
$ws.Range('G36').Value = '23'
$ws.Range('G37').Value = '23'
$ws.Range('G38').Value = '23'
$ws.Range('G39').Value = '23'
$ws.Range('D40').Value = '0.03919'
$ws.Range('G40').Value = '23'
$ws.Range('D41').Value = '0.006819'
$ws.Range('G41').Value = '23'
$ws.Range('D42').Value = '0.1064'
$ws.Range('G42').Value = '23'
$ws.Range('D43').Value = '0.003202'
$ws.Range('G43').Value = '23'
$ws.Range('D44').Value = '0.007469'
$ws.Range('G44').Value = '23'
$ws.Range('D45').Value = '0.00005597'
$ws.Range('G45').Value = '23'
$ws.Range('G46').Value = '23'
$ws.Range('D47').Value = '0.4502'
$ws.Range('G47').Value = '23'
$ws.Range('D48').Value = '0.1964'
$ws.Range('G48').Value = '23'
$ws.Range('G49').Value = '23'
$ws.Range('G50').Value = '23'
$ws.Range('G51').Value = '23'

# --- Restore default style on the numeric-text cells ---
foreach ($c in $numCells) {
    $ws.Range($c).Style = "Normal"
}
